$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "subject" index column (A) was off by one for rows 17-30 (values 16..29);
# correct it so each row's subject number matches its row number (17..30).
for ($r = 17; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = $r
}

# Reflect the saved view state: scrolled down with A30 as the active cell.
$ws.Range("A30").Select()
